$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new result rows (11 and 12) mirroring the existing result-row layout,
# fixing the AVR (average transaction/value rate?) pending-time figures.

$ws.Cells.Item(11, 1).Value = 2
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = 10
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 3
$ws.Cells.Item(11, 7).Value = 50
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 200
$ws.Cells.Item(11, 12).Value = $false
$ws.Cells.Item(11, 13).Value = 16384
$ws.Cells.Item(11, 14).Value = 40
$ws.Cells.Item(11, 15).Value = "<-parameter / result->"
$ws.Cells.Item(11, 16).Value = 8
$ws.Cells.Item(11, 17).Value = 0.2571566200256348
$ws.Cells.Item(11, 17).NumberFormat = "0.000000"
$ws.Cells.Item(11, 18).Value = 111.4790219664574
$ws.Cells.Item(11, 18).NumberFormat = "0.000000"
$ws.Cells.Item(11, 19).Value = 32.14457750320435
$ws.Cells.Item(11, 19).NumberFormat = "0.000000"
$ws.Cells.Item(11, 20).Value = 0.2571566200256348
$ws.Cells.Item(11, 20).NumberFormat = "0.000000"
$ws.Cells.Item(11, 21).Value = 0.00058746337890625
$ws.Cells.Item(11, 21).NumberFormat = "0.000000"
$ws.Cells.Item(11, 22).Value = 0.08609374999999998
$ws.Cells.Item(11, 22).NumberFormat = "0.000000"
$ws.Cells.Item(11, 23).Value = 0.04304687499999999
$ws.Cells.Item(11, 23).NumberFormat = "0.000000"
$ws.Cells.Item(11, 24).Value = 0
$ws.Cells.Item(11, 24).NumberFormat = "0.000000"
$ws.Cells.Item(11, 25).Value = 70.0704345703125
$ws.Cells.Item(11, 25).NumberFormat = "0.000"
$ws.Cells.Item(11, 26).Value = 8676.8
$ws.Cells.Item(11, 26).NumberFormat = "0.00"
$ws.Cells.Item(11, 27).Value = 1856
$ws.Cells.Item(11, 27).NumberFormat = "0.00"
$ws.Cells.Item(12, 1).Value = 2
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = 10
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = 3
$ws.Cells.Item(12, 7).Value = 100
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 200
$ws.Cells.Item(12, 12).Value = $false
$ws.Cells.Item(12, 13).Value = 32768
$ws.Cells.Item(12, 14).Value = 95
$ws.Cells.Item(12, 15).Value = "<-parameter / result->"
$ws.Cells.Item(12, 16).Value = 15
$ws.Cells.Item(12, 17).Value = 0.5498097991943359
$ws.Cells.Item(12, 17).NumberFormat = "0.000000"
$ws.Cells.Item(12, 18).Value = 245.275257277907
$ws.Cells.Item(12, 18).NumberFormat = "0.000000"
$ws.Cells.Item(12, 19).Value = 36.65398661295573
$ws.Cells.Item(12, 19).NumberFormat = "0.000000"
$ws.Cells.Item(12, 20).Value = 0.5498097991943359
$ws.Cells.Item(12, 20).NumberFormat = "0.000000"
$ws.Cells.Item(12, 21).Value = 0.0004422664642333984
$ws.Cells.Item(12, 21).NumberFormat = "0.000000"
$ws.Cells.Item(12, 22).Value = 0.1540625
$ws.Cells.Item(12, 22).NumberFormat = "0.000000"
$ws.Cells.Item(12, 23).Value = 0.07703125
$ws.Cells.Item(12, 23).NumberFormat = "0.000000"
$ws.Cells.Item(12, 24).Value = 0
$ws.Cells.Item(12, 24).NumberFormat = "0.000000"
$ws.Cells.Item(12, 25).Value = 140.8762097358704
$ws.Cells.Item(12, 25).NumberFormat = "0.000"
$ws.Cells.Item(12, 26).Value = 17585.6
$ws.Cells.Item(12, 26).NumberFormat = "0.00"
$ws.Cells.Item(12, 27).Value = 3480
$ws.Cells.Item(12, 27).NumberFormat = "0.00"

Write-Host "rows 11-12 written"
